$d = $word.ActiveDocument

# Update the title paragraph (date)
$d.Paragraphs.Item(1).Range.Text = "2025-09-15 Monday"

# Update each cell of the table with the new arithmetic expression
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "68+13=81"
$t.Cell(1, 2).Range.Text = "16+29=45"
$t.Cell(1, 3).Range.Text = "19+43=62"
$t.Cell(1, 4).Range.Text = "29+16=45"
$t.Cell(1, 5).Range.Text = "8+23=31"
$t.Cell(2, 1).Range.Text = "64-5=59"
$t.Cell(2, 2).Range.Text = "30-26=4"
$t.Cell(2, 3).Range.Text = "5+16=21"
$t.Cell(2, 4).Range.Text = "33+9=42"
$t.Cell(2, 5).Range.Text = "59+25=84"
$t.Cell(3, 1).Range.Text = "15+58=73"
$t.Cell(3, 2).Range.Text = "85-58=27"
$t.Cell(3, 3).Range.Text = "48+36=84"
$t.Cell(3, 4).Range.Text = "38+5=43"
$t.Cell(3, 5).Range.Text = "58+7=65"
$t.Cell(4, 1).Range.Text = "37+59=96"
$t.Cell(4, 2).Range.Text = "25+18=43"
$t.Cell(4, 3).Range.Text = "72-13=59"
$t.Cell(4, 4).Range.Text = "86-49=37"
$t.Cell(4, 5).Range.Text = "86-29=57"
$t.Cell(5, 1).Range.Text = "40-36=4"
$t.Cell(5, 2).Range.Text = "38+55=93"
$t.Cell(5, 3).Range.Text = "9+56=65"
$t.Cell(5, 4).Range.Text = "90-9=81"
$t.Cell(5, 5).Range.Text = "77+6=83"
$t.Cell(6, 1).Range.Text = "41-13=28"
$t.Cell(6, 2).Range.Text = "62-8=54"
$t.Cell(6, 3).Range.Text = "92-53=39"
$t.Cell(6, 4).Range.Text = "97-88=9"
$t.Cell(6, 5).Range.Text = "65-39=26"
$t.Cell(7, 1).Range.Text = "87+9=96"
$t.Cell(7, 2).Range.Text = "90-19=71"
$t.Cell(7, 3).Range.Text = "18+4=22"
$t.Cell(7, 4).Range.Text = "73-65=8"
$t.Cell(7, 5).Range.Text = "67+16=83"
$t.Cell(8, 1).Range.Text = "94-47=47"
$t.Cell(8, 2).Range.Text = "46+45=91"
$t.Cell(8, 3).Range.Text = "28+63=91"
$t.Cell(8, 4).Range.Text = "67+27=94"
$t.Cell(8, 5).Range.Text = "27+45=72"
$t.Cell(9, 1).Range.Text = "14+79=93"
$t.Cell(9, 2).Range.Text = "54+29=83"
$t.Cell(9, 3).Range.Text = "58+23=81"
$t.Cell(9, 4).Range.Text = "82-35=47"
$t.Cell(9, 5).Range.Text = "39+19=58"
$t.Cell(10, 1).Range.Text = "8+17=25"
$t.Cell(10, 2).Range.Text = "65-38=27"
$t.Cell(10, 3).Range.Text = "84-65=19"
$t.Cell(10, 4).Range.Text = "27+18=45"
$t.Cell(10, 5).Range.Text = "17+68=85"
$t.Cell(11, 1).Range.Text = "25+66=91"
$t.Cell(11, 2).Range.Text = "55-8=47"
$t.Cell(11, 3).Range.Text = "93-76=17"
$t.Cell(11, 4).Range.Text = "7+68=75"
$t.Cell(11, 5).Range.Text = "46-29=17"
$t.Cell(12, 1).Range.Text = "45+29=74"
$t.Cell(12, 2).Range.Text = "75-16=59"
$t.Cell(12, 3).Range.Text = "28+26=54"
$t.Cell(12, 4).Range.Text = "46+49=95"
$t.Cell(12, 5).Range.Text = "51-42=9"
$t.Cell(13, 1).Range.Text = "81-2=79"
$t.Cell(13, 2).Range.Text = "58+15=73"
$t.Cell(13, 3).Range.Text = "55-17=38"
$t.Cell(13, 4).Range.Text = "12+9=21"
$t.Cell(13, 5).Range.Text = "93-58=35"
$t.Cell(14, 1).Range.Text = "28+48=76"
$t.Cell(14, 2).Range.Text = "95-37=58"
$t.Cell(14, 3).Range.Text = "70-67=3"
$t.Cell(14, 4).Range.Text = "17+35=52"
$t.Cell(14, 5).Range.Text = "45+17=62"
$t.Cell(15, 1).Range.Text = "39+53=92"
$t.Cell(15, 2).Range.Text = "6+86=92"
$t.Cell(15, 3).Range.Text = "67-8=59"
$t.Cell(15, 4).Range.Text = "74-17=57"
$t.Cell(15, 5).Range.Text = "87-59=28"
$t.Cell(16, 1).Range.Text = "35+18=53"
$t.Cell(16, 2).Range.Text = "39+53=92"
$t.Cell(16, 3).Range.Text = "90-77=13"
$t.Cell(16, 4).Range.Text = "42-19=23"
$t.Cell(16, 5).Range.Text = "19+73=92"
$t.Cell(17, 1).Range.Text = "94-86=8"
$t.Cell(17, 2).Range.Text = "70-15=55"
$t.Cell(17, 3).Range.Text = "25-16=9"
$t.Cell(17, 4).Range.Text = "48+49=97"
$t.Cell(17, 5).Range.Text = "84-19=65"
$t.Cell(18, 1).Range.Text = "15+77=92"
$t.Cell(18, 2).Range.Text = "57+8=65"
$t.Cell(18, 3).Range.Text = "91-76=15"
$t.Cell(18, 4).Range.Text = "92-26=66"
$t.Cell(18, 5).Range.Text = "83-55=28"
$t.Cell(19, 1).Range.Text = "85-28=57"
$t.Cell(19, 2).Range.Text = "51-47=4"
$t.Cell(19, 3).Range.Text = "80-78=2"
$t.Cell(19, 4).Range.Text = "25+48=73"
$t.Cell(19, 5).Range.Text = "28+4=32"
$t.Cell(20, 1).Range.Text = "15+17=32"
$t.Cell(20, 2).Range.Text = "54+9=63"
$t.Cell(20, 3).Range.Text = "2+79=81"
$t.Cell(20, 4).Range.Text = "58+38=96"
$t.Cell(20, 5).Range.Text = "34+47=81"
